$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.683891415596008
$ws.Range("B1").Value = 2.119223833084106
$ws.Range("C1").Value = 2.02674126625061
$ws.Range("D1").Value = 1.617934107780457
$ws.Range("E1").Value = 1.48410964012146
